$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: the existing "总计" sheet (sheetId 6) becomes "2022-Q1".
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Clear everything on the sheet first so no stale cells survive.
$q1.Cells.ClearContents()

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy the header style (bold + border, centered) from column B onto the
# newly introduced header cells D1:H1 (B1/C1/D1 already carried it from
# the original sheet, but D1 must be re-stamped since its text changed
# and E1:H1 are brand new).
$q1.Range("B1").Copy($q1.Range("D1"))
$q1.Range("D1").Value = "基金规模"
$q1.Range("B1").Copy($q1.Range("E1"))
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("B1").Copy($q1.Range("F1"))
$q1.Range("F1").Value = "仓位占比"
$q1.Range("B1").Copy($q1.Range("G1"))
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("B1").Copy($q1.Range("H1"))
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "001675"
$q1.Range("C2").Value = "江信同福灵活配置混合A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "1.02"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "92.43"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "4.05"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0413"
$q1.Range("H2").Value = 10

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "001676"
$q1.Range("C3").Value = "江信同福灵活配置混合C"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "0.42"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "92.43"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "4.05"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.0170"
$q1.Range("H3").Value = 10

# Re-stamp the A-column style (bold/border/centered) onto A2/A3, which
# already carried it from the old sheet, to be safe after ClearContents.
$q1.Range("B1").Copy($q1.Range("A2"))
$q1.Range("A2").Value = 0
$q1.Range("B1").Copy($q1.Range("A3"))
$q1.Range("A3").Value = 1

# ------------------------------------------------------------------
# Step 2: add a brand-new "总计" sheet after "2022-Q1", holding the
# updated roll-up table (2022-Q1 row added, rest shifted down).
# ------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add()
$newTotal.Name = "总计"
$q1after = $wb.Worksheets.Item("2022-Q1")
$newTotal.Move($null, $q1after)

$total = $wb.Worksheets.Item("总计")
$total.Cells.ClearContents()

# Header row (reuse style from the "2022-Q1" sheet's header cell so the
# bold+border formatting matches exactly).
$q1after.Range("B1").Copy($total.Range("B1"))
$total.Range("B1").Value = "日期"
$q1after.Range("B1").Copy($total.Range("C1"))
$total.Range("C1").Value = "持有数量(只)"
$q1after.Range("B1").Copy($total.Range("D1"))
$total.Range("D1").Value = "持有市值(亿元)"

$rows = @(
    @{A=0; B="2022-Q1"; C=2; D=0.06},
    @{A=1; B="2021-Q4"; C=5; D=2.01},
    @{A=2; B="2021-Q3"; C=5; D=0.79},
    @{A=3; B="2021-Q2"; C=3; D=0.13},
    @{A=4; B="2021-Q1"; C=2; D=0.05},
    @{A=5; B="2020-Q4"; C=2; D=0.11}
)

$r = 2
foreach ($row in $rows) {
    $q1after.Range("B1").Copy($total.Range("A$r"))
    $total.Range("A$r").Value = $row.A
    $total.Range("B$r").Value = $row.B
    $total.Range("C$r").Value = $row.C
    $total.Range("D$r").Value = $row.D
    $r = $r + 1
}
